$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aValues = @(5490,5460,5430,5400,5380,5350,5340,5350,5360,5370,5380,5400,5420,5460,5520,5610,5720,5860,6020,6190,6370,6550,6740,6910,7080,7220,7350,7450,7530,7590,7620,7610,7570,7520,7470,7400,7340,7280,7230,7190,7160,7130,7110,7100,7110,7120,7130,7140,7150,7170,7200,7250,7310,7390,7480,7580,7670,7750,7810,7840,7850,7840,7830,7800,7760,7700,7650,7570,7500,7410,7300,7170,7030,6870,6720,6550,6430,6270,6140,6020,5900,5830,5770,5720,5680)
$bValues = @(45972,45972.01041666666,45972.02083333334,45972.03125,45972.04166666666,45972.05208333334,45972.0625,45972.09375,45972.10416666666,45972.11458333334,45972.125,45972.13541666666,45972.14583333334,45972.15625,45972.16666666666,45972.17708333334,45972.1875,45972.19791666666,45972.20833333334,45972.21875,45972.22916666666,45972.23958333334,45972.25,45972.26041666666,45972.27083333334,45972.28125,45972.29166666666,45972.30208333334,45972.3125,45972.32291666666,45972.33333333334,45972.35416666666,45972.36458333334,45972.375,45972.38541666666,45972.39583333334,45972.40625,45972.41666666666,45972.42708333334,45972.4375,45972.44791666666,45972.45833333334,45972.46875,45972.47916666666,45972.53125,45972.55208333334,45972.5625,45972.57291666666,45972.58333333334,45972.59375,45972.60416666666,45972.61458333334,45972.625,45972.63541666666,45972.64583333334,45972.65625,45972.66666666666,45972.67708333334,45972.6875,45972.69791666666,45972.70833333334,45972.72916666666,45972.73958333334,45972.77083333334,45972.78125,45972.79166666666,45972.80208333334,45972.8125,45972.82291666666,45972.83333333334,45972.84375,45972.85416666666,45972.86458333334,45972.875,45972.88541666666,45972.89583333334,45972.90625,45972.91666666666,45972.92708333334,45972.9375,45972.94791666666,45972.95833333334,45972.96875,45972.97916666666,45972.98958333334)

for ($i = 0; $i -lt $aValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

# Remove now-obsolete trailing rows (previously rows 87-97)
$ws.Range("A87:A97").EntireRow.Delete() | Out-Null
